# "Refactored and created general menu" - update the Sprint2 burndown
# worksheet: the team's availability for 23-27 Jan was revised down, a
# backlog item's description/assignee got corrected, and the selection
# was left on G36 (where the edit ended).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint2")

# Row 36 (Design - HighScore): available hours for G..K trimmed from 1 to 0.5
$ws.Range("G36:K36").Value = 0.5

# Row 17/37 (StartMenu - Magnus/Tobias/Stoffe): Mon estimate 1.5 -> 1,
# Tue availability 0.5 -> 0
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 0

# Row 4/38 (CountDown - Philip): availability for Tue..Sat dropped to 0
$ws.Range("G38:K38").Value = 0

# Row 23/43: backlog item description corrected, and resource reassigned
$ws.Range("D43").Value = "Skriva in sitt namn"
$ws.Range("E43").Value = "Tobias/Stoffe"
$ws.Range("G43:K43").Value = 0

# Leave the active selection where the edits were made
[void]$ws.Range("G36").Select()

# Recalculate dependent totals (row 60 sums, and the burndown rows 69:73)
[void]$wb.Application.CalculateFull()
